# Applies the data update from the commit diff to the "ResumoInscricoes" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 114
$ws.Range("F2").Value = 84
$ws.Range("H2").Value = 90
$ws.Range("E5").Value = 172
$ws.Range("F5").Value = 119
$ws.Range("H5").Value = 130
$ws.Range("E6").Value = 55
$ws.Range("E7").Value = 46
$ws.Range("E10").Value = 745
$ws.Range("F10").Value = 427
$ws.Range("H10").Value = 522
$ws.Range("E11").Value = 499
$ws.Range("F11").Value = 292
$ws.Range("H11").Value = 357
$ws.Range("E12").Value = 766
$ws.Range("F12").Value = 469
$ws.Range("H12").Value = 555
$ws.Range("F13").Value = 100
$ws.Range("H13").Value = 134
$ws.Range("E14").Value = 152
$ws.Range("F14").Value = 87
$ws.Range("H14").Value = 121
$ws.Range("E15").Value = 213
$ws.Range("F15").Value = 101
$ws.Range("H15").Value = 152
$ws.Range("E16").Value = 241
$ws.Range("F16").Value = 140
$ws.Range("H16").Value = 188
$ws.Range("F20").Value = 50
$ws.Range("H20").Value = 87
$ws.Range("F21").Value = 91
$ws.Range("H21").Value = 122
$ws.Range("E22").Value = 205
$ws.Range("F22").Value = 116
$ws.Range("H22").Value = 158
$ws.Range("E23").Value = 236
$ws.Range("F23").Value = 122
$ws.Range("H23").Value = 174
$ws.Range("E24").Value = 290
$ws.Range("F24").Value = 167
$ws.Range("H24").Value = 197
$ws.Range("E25").Value = 353
$ws.Range("F25").Value = 194
$ws.Range("H25").Value = 254
$ws.Range("F26").Value = 134
$ws.Range("H26").Value = 159
$ws.Range("E27").Value = 406
$ws.Range("F27").Value = 225
$ws.Range("H27").Value = 307
$ws.Range("E28").Value = 236
$ws.Range("F28").Value = 116
$ws.Range("H28").Value = 168
$ws.Range("E29").Value = 202
$ws.Range("F29").Value = 119
$ws.Range("H29").Value = 160
$ws.Range("E30").Value = 271
$ws.Range("F30").Value = 170
$ws.Range("H30").Value = 223
$ws.Range("F31").Value = 40
$ws.Range("H31").Value = 67
$ws.Range("F32").Value = 145
$ws.Range("H32").Value = 183
$ws.Range("E33").Value = 352
$ws.Range("F33").Value = 185
$ws.Range("H33").Value = 276
$ws.Range("E34").Value = 269
$ws.Range("F34").Value = 188
$ws.Range("H34").Value = 226
$ws.Range("E35").Value = 192
$ws.Range("F35").Value = 134
$ws.Range("H35").Value = 161
$ws.Range("F36").Value = 60
$ws.Range("H36").Value = 70
$ws.Range("F37").Value = 113
$ws.Range("H37").Value = 149
$ws.Range("E38").Value = 110
$ws.Range("F38").Value = 70
$ws.Range("H38").Value = 87
$ws.Range("E39").Value = 210
$ws.Range("F39").Value = 108
$ws.Range("H39").Value = 159
$ws.Range("E40").Value = 320
$ws.Range("F40").Value = 174
$ws.Range("H40").Value = 254
$ws.Range("E41").Value = 459
$ws.Range("F41").Value = 235
$ws.Range("H41").Value = 327
$ws.Range("E42").Value = 489
$ws.Range("F42").Value = 283
$ws.Range("H42").Value = 344
$ws.Range("E43").Value = 152
$ws.Range("F43").Value = 88
$ws.Range("H43").Value = 115
$ws.Range("E44").Value = 394
$ws.Range("F44").Value = 213
$ws.Range("H44").Value = 281
$ws.Range("E45").Value = 191
$ws.Range("F45").Value = 112
$ws.Range("H45").Value = 151
$ws.Range("E46").Value = 410
$ws.Range("F46").Value = 243
$ws.Range("H46").Value = 307
$ws.Range("E47").Value = 569
$ws.Range("F47").Value = 325
$ws.Range("H47").Value = 417
$ws.Range("E48").Value = 290
$ws.Range("F48").Value = 140
$ws.Range("H48").Value = 184
$ws.Range("F49").Value = 175
$ws.Range("H49").Value = 262
$ws.Range("F50").Value = 164
$ws.Range("H50").Value = 237
$ws.Range("E51").Value = 268
$ws.Range("F51").Value = 139
$ws.Range("H51").Value = 213
